# "Small changes to test data"
#
# registerValidUserSheet: bump the sample username/password test values
# registerInvalidUserSheet: drop the now-redundant "Qwerty)123" duplicate row

$wb = $excel.ActiveWorkbook

# --- registerValidUserSheet --------------------------------------------
# A2/A3 feed the CONCAT() formulas in D2/D3, so those recalc automatically.
$ws1 = $wb.Worksheets.Item("registerValidUserSheet")
$ws1.Range("A2").Value = "testuserAvengers35"
$ws1.Range("A3").Value = "Qwerty+1234586"

# --- registerInvalidUserSheet -------------------------------------------
# Remove row 3 (username "Qwerty)123") entirely; remaining rows shift up.
$ws2 = $wb.Worksheets.Item("registerInvalidUserSheet")
$ws2.Rows.Item(3).Delete()
$ws2.Range("A3:XFD3").Select()

# Restore the workbook's original active sheet/tab.
$ws3 = $wb.Worksheets.Item("practiceQuestionSheet")
$ws3.Activate()
